$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

# Remove the four stray "it_*" rows (rows 2-5) so the UK sheet starts
# directly with the uk_*_std.txt rows; remaining rows shift up and the
# sheet's used range shrinks from A1:D21 to A1:D17.
$ws.Rows("2:5").Delete()
